# T3 status report update - "edits to the reports"
# Updates the title and the Key Accomplishments section of the status report
# with the team's real names/roles and a fuller description of work done.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Key Accomplishments table: flesh out names with roles and describe the work
# in more detail.
$ws.Range("B11").Value = "reviewed aspects of order books and went over more requirements for the project "
$ws.Range("B12").Value = "reviewed aspects of order books and went over more requirements for the project "
$ws.Range("B10").Value = " reviewed aspects of order books and went over more requirements for the project "

# Re-title the report for the real team/project (was a generic class-project title)
$ws.Range("A1").Value = "Status Report for Team 3 - Orderbook  Project"

$ws.Range("A10").Value = "Zeenat (Lead)"
$ws.Range("A11").Value = "Buchi (frontend)"
$ws.Range("A12").Value = "Anmol(backend)"

# Move the active selection to B5 (contact cell) as in the saved workbook
$ws.Range("B5").Select()
